# Update the "repaymentstrategy" value on the ProductLoanInput sheet
# from "RBI (India)" to "Overdue/Due Fee/Int,Principal"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the active selection moving to B17, as captured in the saved view state
$ws.Activate()
$ws.Range("B17").Select()
